$wb = $excel.ActiveWorkbook

# --- Sheet: compare_models (column I updates) ---
$ws1 = $wb.Worksheets.Item("compare_models")
$ws1.Range("I2").Value = 0.074
$ws1.Range("I3").Value = 0.056
$ws1.Range("I4").Value = 0.08599999999999999
$ws1.Range("I6").Value = 0.046
$ws1.Range("I7").Value = 1.05
$ws1.Range("I8").Value = 0.018
$ws1.Range("I10").Value = 0.022
$ws1.Range("I11").Value = 0.024
$ws1.Range("I12").Value = 0.524
$ws1.Range("I13").Value = 0.024
$ws1.Range("I14").Value = 0.018
$ws1.Range("I15").Value = 0.016
$ws1.Range("I16").Value = 0.02
$ws1.Range("I17").Value = 0.018

# --- Sheet: pred_final (row 2, columns C:H updates) ---
$ws2 = $wb.Worksheets.Item("pred_final")
$ws2.Range("C2").Value = 1.1033
$ws2.Range("D2").Value = 2.9955
$ws2.Range("E2").Value = 1.7308
$ws2.Range("F2").Value = 0.9955000000000001
$ws2.Range("G2").Value = 0.0342
$ws2.Range("H2").Value = 0.0192
